# Update TPM-derived NATMI ligand-receptor metrics (Adam15-Itgb3) with
# refreshed expression values (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.749936
$ws.Range("H2").Value = 95.249808
$ws.Range("I2").Value = 0.5302851438878331
$ws.Range("J2").Value = 0.5302851438878331
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 9.213492761216001
$ws.Range("R2").Value = 82.921434850944
$ws.Range("S2").Value = 0.01818554345883485
$ws.Range("T2").Value = 0.01818554345883485
$ws.Range("G3").Value = 31.749936
$ws.Range("H3").Value = 95.249808
$ws.Range("I3").Value = 0.5302851438878331
$ws.Range("J3").Value = 0.5302851438878331
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 225.753178604224
$ws.Range("R3").Value = 2031.778607438016
$ws.Range("S3").Value = 0.4455904342551828
$ws.Range("T3").Value = 0.4455904342551828
$ws.Range("G4").Value = 31.749936
$ws.Range("H4").Value = 95.249808
$ws.Range("I4").Value = 0.5302851438878331
$ws.Range("J4").Value = 0.5302851438878331
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 33.69609066036801
$ws.Range("R4").Value = 303.264815943312
$ws.Range("S4").Value = 0.06650916617381553
$ws.Range("T4").Value = 0.06650916617381551
$ws.Range("I5").Value = 0.3451699599880819
$ws.Range("J5").Value = 0.3451699599880819
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 5.997190312410666
$ws.Range("R5").Value = 53.97471281169599
$ws.Range("S5").Value = 0.01183722263464973
$ws.Range("T5").Value = 0.01183722263464973
$ws.Range("I6").Value = 0.3451699599880819
$ws.Range("J6").Value = 0.3451699599880819
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.2900409980096793
$ws.Range("T6").Value = 0.2900409980096793
$ws.Range("I7").Value = 0.3451699599880819
$ws.Range("J7").Value = 0.3451699599880819
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 21.93325307912867
$ws.Range("S7").Value = 0.04329173934375294
$ws.Range("T7").Value = 0.04329173934375292
$ws.Range("I8").Value = 0.1245448961240849
$ws.Range("J8").Value = 0.1245448961240849
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 2.163917869681778
$ws.Range("R8").Value = 19.475260827136
$ws.Range("S8").Value = 0.004271129687766054
$ws.Range("T8").Value = 0.004271129687766054
$ws.Range("I9").Value = 0.1245448961240849
$ws.Range("J9").Value = 0.1245448961240849
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.1046531568682532
$ws.Range("T9").Value = 0.1046531568682532
$ws.Range("I10").Value = 0.1245448961240849
$ws.Range("J10").Value = 0.1245448961240849
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 7.913999023836444
$ws.Range("R10").Value = 71.22599121452801
$ws.Range("S10").Value = 0.01562060956806565
$ws.Range("T10").Value = 0.01562060956806564
